# Apply cryptos-list price/volume update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.528.00'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '2.524.57'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.65'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.88'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('E7').Value = '  -1.28%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -2.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.15'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.82%  '
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.17'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.29%  '
$ws.Range('D14').Value = '2.912.59'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.28'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.64%  '
$ws.Range('D16').Value = '2.508.45'
$ws.Range('E16').Value = '  -1.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.807'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.35%  '
$ws.Range('D18').Value = '42.533.68'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('E19').Value = '  -3.01%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.16'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.37%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0938'
$ws.Range('E21').Value = '  -1.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.96'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '241.27'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.32%  '
$ws.Range('E24').Value = '  -2.56%  '
$ws.Range('E25').Value = '  -3.21%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.45'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.53%  '
$ws.Range('E28').Value = '  -4.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.97'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.74'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -6.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.86'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.78%  '
$ws.Range('E32').Value = '  -0.93%  '
$ws.Range('E33').Value = '  -3.13%  '
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0782'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.80%  '
$ws.Range('E36').Value = '  -1.81%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.97'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.64%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.61'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.108'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.10%  '
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.17'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.15'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.08%  '
$ws.Range('E44').Value = '  -1.15%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.004.18'
$ws.Range('E45').Value = '  +1.86%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.21'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.56%  '
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').Value = '2.760.97'
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '79.24'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.188'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.45'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.28%  '
